# Applies the "Reporte de entrega del 40%" update to the
# "Lista de Tareas de la 4ta Iteracion" workbook (sheet "Casos de Uso"):
#   - Updates the Estatus (status) column for several tasks.
#   - Registers consumed hours for Dia 10 (column AI) on rows 13, 14 and 16,
#     which cascades through the dependent Rest./Total formulas.
#   - Updates the active-cell selection left behind by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Estatus column (F) updates -----------------------------------------
$ws.Range("F11").Value = "Hecho"
$ws.Range("F12").Value = "Hecho"
$ws.Range("F13").Value = "En proceso"
$ws.Range("F14").Value = "Hecho"
$ws.Range("F16").Value = "Hecho"
$ws.Range("F17").Value = "Hecho"

# --- Consumed hours for Dia 10 (column AI) -------------------------------
$ws.Range("AI13").Value = 1
$ws.Range("AI14").Value = 2
$ws.Range("AI16").Value = 3

# Force a full recalculation so the dependent "Rest."/"Total" formulas
# (AJ:BA) pick up the new consumed-hours figures.
$excel.CalculateFullRebuild()

# --- Re-touch the day-total merged header cells --------------------------
# (mirrors the re-ordering of the mergeCells list seen when Excel re-saves
# the sheet after the edit)
$touchedMerges = @("AZ4:BA4", "AO4:AP4", "AR4:AS4", "AU4:AV4", "AX4:AY4")
foreach ($r in $touchedMerges) {
    $ws.Range($r).UnMerge()
}
foreach ($r in $touchedMerges) {
    $ws.Range($r).Merge()
}

# --- Final selection state -------------------------------------------------
$ws.Activate()
$ws.Range("AI15").Select()
